$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are non-numeric-looking strings (contain significance
# stars, brackets, etc.) can be assigned directly; Excel will store them as text.
$directValues = @{
    "L4" = "0.07**"
    "L5" = "-0.05*"
    "M5" = "0.15***"
    "L6" = "-0.06*"
    "M6" = "0.09***"
    "L7" = "-0.15***"
    "L9" = "-0.06*"
    "M11" = "0.07**"
    "L12" = "[0.08]"
    "M12" = "-0.44***"
    "P12" = "0.07**"
    "Q12" = "0.06**"
    "M13" = "[0.29]"
    "M14" = "0.36*"
}
foreach ($addr in $directValues.Keys) {
    $ws.Range($addr).Value = $directValues[$addr]
}

# Cells whose new values look like plain numbers must be forced to remain text
# (matching the workbook's original inlineStr/text cell type) by temporarily
# applying a text number format, then restoring the original cell style so no
# visible formatting change is introduced.
$textValues = @{
    "L2" = "0.03"
    "L3" = "0.03"
    "L8" = "0.05"
    "M8" = "-0.04"
    "M10" = "0.03"
    "L11" = "0.03"
    "B12" = "0.03"
    "C12" = "-0.02"
    "D12" = "0.04"
    "E12" = "-0.20"
    "F12" = "0.09"
    "G12" = "-0.16"
    "H12" = "0.20"
    "J12" = "0.03"
    "N12" = "-0.01"
    "R12" = "-0.04"
    "D13" = "-0.08"
    "F13" = "0.05"
    "G13" = "0.02"
    "H13" = "-0.05"
    "L13" = "-0.13"
    "L14" = "-0.26"
    "L15" = "-0.01"
    "M15" = "0.21"
    "L16" = "-0.10"
    "M16" = "0.23"
    "M17" = "-0.02"
    "L18" = "0.12"
    "M18" = "-0.24"
}
foreach ($addr in $textValues.Keys) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $textValues[$addr]
    $r.Style = $origStyle
}
